$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.644.24'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -4.29%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.337.75'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.14%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.67'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.30%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '181.07'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.21%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.633'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.93%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  -3.35%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.65'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.49%  '
$ws.Range("E11").Value = '  -3.28%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.918.06'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.20%  '
$ws.Range("E13").Value = '  -0.48%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.04'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.49%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '66.745.00'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.14%  '
$ws.Range("E16").Value = '  -2.83%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.339.36'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.83%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '437.03'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.75%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.69'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.20%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.54'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.06%  '
$ws.Range("E21").Value = '  -2.37%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.45'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.07%  '
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("E24").Value = '  -0.65%  '
$ws.Range("E25").Value = '  -3.95%  '
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.03'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.12%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.79'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.50%  '
$ws.Range("E31").Value = '  +0.07%  '
$ws.Range("E32").Value = '  -5.60%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.78'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.23%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.22'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.73%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '164.63'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("E36").Value = '  -6.21%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '27.35'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.38%  '
$ws.Range("E38").Value = '  -8.08%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.838.41'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.60%  '
$ws.Range("E40").Value = '  -1.73%  '
$ws.Range("E41").Value = '  -3.82%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.21'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.88%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.21'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.14%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0668'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.92%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '24.41'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.11%  '
$ws.Range("E46").Value = '  -7.16%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '321.33'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.59%  '
$ws.Range("E48").Value = '  -3.57%  '
$ws.Range("E49").Value = '  +0.95%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.976'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.40%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.16'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.63%  '
